$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 75 (pushes existing rows 75-113 down to 76-114)
$ws.Rows.Item(75).Insert()

$ws.Range("A75").Value = 4
$ws.Range("B75").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C75").Value = "Los Lagos"
$ws.Range("D75").Value = 44488
$ws.Range("E75").Value = 10
$ws.Range("F75").Value = 100112009
$ws.Range("G75").Value = "Acelga"
$ws.Range("H75").Value = "Sin especificar"
$ws.Range("I75").Value = "Primera"
$ws.Range("J75").Value = 200
$ws.Range("K75").Value = 3500
$ws.Range("L75").Value = 3500
$ws.Range("M75").Value = 3500
$ws.Range("N75").Value = "`$/docena de atados (4 kilos)"
$ws.Range("O75").Value = "Región del Maule"
$ws.Range("P75").Value = 875
$ws.Range("Q75").Value = 4
$ws.Range("R75").Value = "Hortaliza"
